$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7173.125
$ws.Range("I32").Value = 7325.154
$ws.Range("J32").Value = 6993.4546
$ws.Range("K32").Value = 7325.154
$ws.Range("L32").Value = 6993.4546
$ws.Range("M32").Value = -6999.154
$ws.Range("N32").Value = -7645.4546
$ws.Range("H43").Value = 4133.778
$ws.Range("J43").Value = 5447.6665
$ws.Range("L43").Value = 5447.6665
$ws.Range("N43").Value = -5585.6665
$ws.Range("H64").Value = 5803.143
$ws.Range("I64").Value = 2684.4
$ws.Range("K64").Value = 2684.4
$ws.Range("M64").Value = -2436.4
$ws.Range("H67").Value = 5803.143
$ws.Range("I67").Value = 2684.4
$ws.Range("K67").Value = 2684.4
$ws.Range("M67").Value = -1826.4
$ws.Range("H138").Value = 1644.1333
$ws.Range("I138").Value = 1079.6316
$ws.Range("J138").Value = 2619.182
$ws.Range("K138").Value = 3238.8948
$ws.Range("L138").Value = 7857.545999999999
$ws.Range("M138").Value = 1901.1052
$ws.Range("N138").Value = -18137.546

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4384.3335
$ws.Range("I32").Value = 4020.64
$ws.Range("K32").Value = 4020.64
$ws.Range("M32").Value = -3733.64
$ws.Range("H74").Value = 30305540
$ws.Range("J74").Value = 58827132
$ws.Range("L74").Value = 58827132
$ws.Range("N74").Value = -58828880
$ws.Range("H77").Value = 30305540
$ws.Range("J77").Value = 58827132
$ws.Range("L77").Value = 294135660
$ws.Range("N77").Value = -294144396
$ws.Range("H122").Value = 13890122
$ws.Range("I122").Value = 1294.4
$ws.Range("K122").Value = 3883.2
$ws.Range("M122").Value = -1433.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 33335224
$ws.Range("I86").Value = 41668240
$ws.Range("J86").Value = 3162.3333
$ws.Range("K86").Value = 41668240
$ws.Range("L86").Value = 3162.3333
$ws.Range("M86").Value = -41667117
$ws.Range("N86").Value = -5408.3333
$ws.Range("H89").Value = 33335224
$ws.Range("I89").Value = 41668240
$ws.Range("J89").Value = 3162.3333
$ws.Range("K89").Value = 208341200
$ws.Range("L89").Value = 15811.6665
$ws.Range("M89").Value = -208335584
$ws.Range("N89").Value = -27043.6665
$ws.Range("H134").Value = 2466.8245
$ws.Range("I134").Value = 1953.196
$ws.Range("J134").Value = 6832.6665
$ws.Range("K134").Value = 5859.588
$ws.Range("L134").Value = 20497.9995
$ws.Range("M134").Value = -3324.588
$ws.Range("N134").Value = -25567.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 70.666664
$ws.Range("J7").Value = 38.28
$ws.Range("L7").Value = 38.28
$ws.Range("N7").Value = -264.28
$ws.Range("H31").Value = 1743.069
$ws.Range("I31").Value = 1338.4762
$ws.Range("J31").Value = 2805.125
$ws.Range("K31").Value = 1338.4762
$ws.Range("L31").Value = 2805.125
$ws.Range("M31").Value = -1043.4762
$ws.Range("N31").Value = -3395.125
$ws.Range("H34").Value = 1743.069
$ws.Range("I34").Value = 1338.4762
$ws.Range("J34").Value = 2805.125
$ws.Range("K34").Value = 1338.4762
$ws.Range("L34").Value = 2805.125
$ws.Range("M34").Value = -1136.4762
$ws.Range("N34").Value = -3209.125
$ws.Range("H58").Value = 33336042
$ws.Range("I58").Value = 20002338
$ws.Range("K58").Value = 20002338
$ws.Range("M58").Value = -20002135
$ws.Range("H62").Value = 8721.5
$ws.Range("I62").Value = 3005
$ws.Range("J62").Value = 10627
$ws.Range("K62").Value = 3005
$ws.Range("L62").Value = 10627
$ws.Range("M62").Value = -2381
$ws.Range("N62").Value = -11875
$ws.Range("H65").Value = 8721.5
$ws.Range("I65").Value = 3005
$ws.Range("J65").Value = 10627
$ws.Range("K65").Value = 15025
$ws.Range("L65").Value = 53135
$ws.Range("M65").Value = -11905
$ws.Range("N65").Value = -59375
$ws.Range("H132").Value = 1945.2609
$ws.Range("I132").Value = 1945.2609
$ws.Range("K132").Value = 5835.7827
$ws.Range("M132").Value = -3305.7827
$ws.Range("H136").Value = 33336042
$ws.Range("I136").Value = 20002338
$ws.Range("K136").Value = 60007014
$ws.Range("M136").Value = -60004464

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 711.7143
$ws.Range("I7").Value = 760.3333
$ws.Range("J7").Value = 420
$ws.Range("K7").Value = 2280.9999
$ws.Range("L7").Value = 1260
$ws.Range("M7").Value = -2168.9999
$ws.Range("N7").Value = -1484
$ws.Range("H68").Value = 1010.7143
$ws.Range("I68").Value = 966.6667
$ws.Range("K68").Value = 2900.0001
$ws.Range("M68").Value = -2089.0001
$ws.Range("H69").Value = 4273.909
$ws.Range("I69").Value = 999
$ws.Range("J69").Value = 4601.4
$ws.Range("K69").Value = 2997
$ws.Range("L69").Value = 13804.2
$ws.Range("M69").Value = -2186
$ws.Range("N69").Value = -15426.2
$ws.Range("H71").Value = 1010.7143
$ws.Range("I71").Value = 966.6667
$ws.Range("K71").Value = 8700.0003
$ws.Range("M71").Value = -4644.0003
$ws.Range("H72").Value = 4273.909
$ws.Range("I72").Value = 999
$ws.Range("J72").Value = 4601.4
$ws.Range("K72").Value = 8991
$ws.Range("L72").Value = 41412.6
$ws.Range("M72").Value = -4935
$ws.Range("N72").Value = -49524.6
$ws.Range("H76").Value = 5666.6665
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 7000
$ws.Range("K76").Value = 9000
$ws.Range("L76").Value = 21000
$ws.Range("M76").Value = -8617
$ws.Range("N76").Value = -21766
$ws.Range("H79").Value = 5666.6665
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 7000
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 21000
$ws.Range("M79").Value = -7674
$ws.Range("N79").Value = -23652
$ws.Range("H80").Value = 18499.8
$ws.Range("I80").Value = 35733.332
$ws.Range("J80").Value = 11114
$ws.Range("K80").Value = 107199.996
$ws.Range("L80").Value = 33342
$ws.Range("M80").Value = -106263.996
$ws.Range("N80").Value = -35214
$ws.Range("H83").Value = 18499.8
$ws.Range("I83").Value = 35733.332
$ws.Range("J83").Value = 11114
$ws.Range("K83").Value = 321599.988
$ws.Range("L83").Value = 100026
$ws.Range("M83").Value = -316919.988
$ws.Range("N83").Value = -109386
$ws.Range("H92").Value = 598.6667
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 598.6667
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1796.0001
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -4292.0001
$ws.Range("H121").Value = 1707.8
$ws.Range("J121").Value = 1897.25
$ws.Range("L121").Value = 5691.75
$ws.Range("N121").Value = -8311.75
$ws.Range("H126").Value = 8501.380999999999
$ws.Range("I126").Value = 1861.3334
$ws.Range("K126").Value = 5584.0002
$ws.Range("M126").Value = -644.0002000000004

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4351.737
$ws.Range("I22").Value = 2598.625
$ws.Range("K22").Value = 2598.625
$ws.Range("M22").Value = -2303.625
$ws.Range("H27").Value = 4351.737
$ws.Range("I27").Value = 2598.625
$ws.Range("K27").Value = 2598.625
$ws.Range("M27").Value = -2491.625
$ws.Range("H93").Value = 958.67566
$ws.Range("I93").Value = 905
$ws.Range("J93").Value = 1037.4
$ws.Range("K93").Value = 905
$ws.Range("L93").Value = 1037.4
$ws.Range("M93").Value = 343
$ws.Range("N93").Value = -3533.4
$ws.Range("H122").Value = 4848.5625
$ws.Range("I122").Value = 3323.4119
$ws.Range("K122").Value = 9970.235700000001
$ws.Range("M122").Value = -7520.235700000001
$ws.Range("H136").Value = 4083905.8
$ws.Range("I136").Value = 2135.4
$ws.Range("K136").Value = 6406.200000000001
$ws.Range("M136").Value = -3856.200000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 729.5
$ws.Range("I107").Value = 837.5
$ws.Range("J107").Value = 621.5
$ws.Range("K107").Value = 2512.5
$ws.Range("L107").Value = 1864.5
$ws.Range("M107").Value = -592.5
$ws.Range("N107").Value = -5704.5
$ws.Range("H113").Value = 1048.9667
$ws.Range("I113").Value = 1084.8334
$ws.Range("J113").Value = 995.1667
$ws.Range("K113").Value = 3254.5002
$ws.Range("L113").Value = 2985.5001
$ws.Range("M113").Value = -1084.5002
$ws.Range("N113").Value = -7325.5001
$ws.Range("H116").Value = 69677.28999999999
$ws.Range("J116").Value = 69677.28999999999
$ws.Range("L116").Value = 69677.28999999999
$ws.Range("N116").Value = -78855.28999999999
$ws.Range("H122").Value = 3597.4375
$ws.Range("I122").Value = 3270
$ws.Range("K122").Value = 9810
$ws.Range("M122").Value = -7360
